# Fruta / hortaliza, semanal
# Insert a new weekly record as row 40, pushing the former rows 40 and 41
# down to rows 41 and 42 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 40 (shifts old row 40 -> 41, old row 41 -> 42)
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new weekly observation
$ws.Range("A40").Value = 5
$ws.Range("B40").Value = "Macroferia Regional de Talca"
$ws.Range("C40").Value = "Maule"
$ws.Range("D40").Value = 44694
$ws.Range("E40").Value = 7
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100107
$ws.Range("H40").Value = "Otros"
$ws.Range("I40").Value = 100107001
$ws.Range("J40").Value = "Caqui"
$ws.Range("K40").Value = "Mankaki"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 180
$ws.Range("N40").Value = 12000
$ws.Range("O40").Value = 12000
$ws.Range("P40").Value = 12000
$ws.Range("Q40").Value = "$/caja 12 kilos granel"
$ws.Range("R40").Value = "Provincia de Curicó"
$ws.Range("S40").Value = 12000
$ws.Range("T40").Value = 1
